$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-30 Wednesday" "2025-05-01 Thursday"

Replace-Text "59×25=1475" "42×79=3318"
Replace-Text "98×66=6468" "77×60=4620"
Replace-Text "57×51=2907" "58×71=4118"
Replace-Text "27×96=2592" "45×42=1890"
Replace-Text "30×17=510" "45×26=1170"

Replace-Text "26×66=1716" "29×14=406"
Replace-Text "43×89=3827" "47×86=4042"
Replace-Text "79×57=4503" "46×35=1610"
Replace-Text "17×49=833" "76×45=3420"
Replace-Text "70×83=5810" "52×16=832"

Replace-Text "81×53=4293" "99×31=3069"
Replace-Text "82×31=2542" "89×16=1424"
Replace-Text "25×14=350" "27×23=621"
Replace-Text "18×66=1188" "55×49=2695"
Replace-Text "38×36=1368" "72×75=5400"

Replace-Text "71×71=5041" "48×16=768"
Replace-Text "87×40=3480" "60×19=1140"
Replace-Text "20×21=420" "44×58=2552"
Replace-Text "27×28=756" "50×20=1000"
Replace-Text "94×58=5452" "78×15=1170"

Replace-Text "31×36=1116" "83×34=2822"
Replace-Text "45×55=2475" "59×66=3894"
Replace-Text "81×81=6561" "72×47=3384"
Replace-Text "52×83=4316" "66×53=3498"
Replace-Text "88×42=3696" "28×89=2492"
